# comentarios v 27 noviembre 2025
# Recompute "Ponderacion_nueva" (column C) on Sheet1 after the
# "041 Alquiler de vivienda" (rent) weight in C44 was adjusted from
# 260 to 237.25. All other category weights are rescaled accordingly
# so that the column continues to add up to 1000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 142.6352067499319
$ws.Range("C3").Value  = 10.68711704379472
$ws.Range("C4").Value  = 6.738949972767558
$ws.Range("C5").Value  = 15.67464814448048
$ws.Range("C6").Value  = 38.9899201284876
$ws.Range("C7").Value  = 12.06165826999728
$ws.Range("C8").Value  = 9.078895941128414
$ws.Range("C9").Value  = 28.09430083801649
$ws.Range("C10").Value = 43.99554765058303
$ws.Range("C11").Value = 9.720138699774713
$ws.Range("C12").Value = 3.948953871958013
$ws.Range("C13").Value = 7.263746193647413
$ws.Range("C14").Value = 1.71601283019335
$ws.Range("C15").Value = 1.666444371549526
$ws.Range("C16").Value = 21.10278776644468
$ws.Range("C17").Value = 19.47961744609215
$ws.Range("C18").Value = 10.03485907211646
$ws.Range("C19").Value = 1.099947701334389
$ws.Range("C20").Value = 29.73635438070952
$ws.Range("C21").Value = 76.96958786175821
$ws.Range("C22").Value = 14.46612191468819
$ws.Range("C23").Value = 0.2336798764637437
$ws.Range("C24").Value = 2.581493854133142
$ws.Range("C25").Value = 27.17374374891689
$ws.Range("C26").Value = 6.945878617582254
$ws.Range("C27").Value = 0.5185018134330206
$ws.Range("C28").Value = 10.82480720669423
$ws.Range("C29").Value = 24.76613290050257
$ws.Range("C30").Value = 10.0600367019038
$ws.Range("C31").Value = 13.02784981308642
$ws.Range("C32").Value = 3.603548263312951
$ws.Range("C33").Value = 1.847408585646028
$ws.Range("C34").Value = 5.059916786324363
$ws.Range("C35").Value = 2.427280871685688
$ws.Range("C36").Value = 86.62520888520288
$ws.Range("C37").Value = 8.199252500433243
$ws.Range("C38").Value = 25.56316224345803
$ws.Range("C39").Value = 5.133089272893819
$ws.Range("C40").Value = 3.557127008392543
$ws.Range("C41").Value = 12.76584510311193
$ws.Range("C42").Value = 0.8324353848439086
$ws.Range("C43").Value = 5.842783712524447
$ws.Range("C44").Value = 237.25
